$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.789.84"
$ws.Range("D3").Value = "3.028.67"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.34"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.41"
$ws.Range("E6").Value = "  -4.62%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "3.028.03"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "3.532.12"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "62.777.30"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "3.029.41"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.16"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.09"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.39"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.51"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.47"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "0.0₃0797"
$ws.Range("E36").Value = "  -6.56%  "
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.28"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  -14.68%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.01"
$ws.Range("E41").Value = "  -3.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "424.46"
$ws.Range("E42").Value = "  -5.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.282"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "2.807.67"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.87"
$ws.Range("E47").Value = "  -9.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.91"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.46"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("E51").Value = "  -1.30%  "
